# Updated cryptos list on Wed Feb  7 03:30:36 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns (and, for two pairs of
# rows whose rank order swapped, the Coin/Link columns too) to match the
# latest scrape. Numeric-looking Price values are written with a leading
# "'" so Excel stores them as text (matching the sheet's existing
# plain-text number formatting, e.g. "1.00" / "96.17") instead of
# re-interpreting them as floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.008.84'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '2.374.04'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = '''96.17'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.497'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = '''34.19'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').Value = '''0.0789'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('E12').Value = '  +2.52%  '
$ws.Range('D13').Value = '''18.22'
$ws.Range('E13').Value = '  -4.23%  '
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').Value = '2.745.48'
$ws.Range('E15').Value = '  +2.21%  '
$ws.Range('D16').Value = '2.376.87'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').Value = '''0.802'
$ws.Range('E17').Value = '  +1.50%  '
$ws.Range('D18').Value = '42.976.66'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').Value = '''12.21'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = '''6.33'
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('D21').Value = '0.0₃0889'
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').Value = '''68.08'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '''235.02'
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '''2.43'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = '''24.89'
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('D28').Value = '''2.37'
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').Value = '''31.57'
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('E33').Value = '  -1.75%  '
$ws.Range('D34').Value = '''0.0732'
$ws.Range('E34').Value = '  +4.31%  '
$ws.Range('E35').Value = '  +5.30%  '
$ws.Range('E36').Value = '  +3.30%  '
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('E38').Value = '  -0.68%  '
$ws.Range('E39').Value = '  +1.49%  '
$ws.Range('D40').Value = '''22.25'
$ws.Range('E40').Value = '  +5.85%  '
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').Value = '''117.83'
$ws.Range('E42').Value = '  -29.16%  '
$ws.Range('D43').Value = '1.943.67'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').Value = '''0.0281'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('D46').Value = '''2.74'
$ws.Range('E46').Value = '  -1.09%  '
$ws.Range('D47').Value = '''9.18'
$ws.Range('E47').Value = '  -10.43%  '
$ws.Range('D48').Value = '2.603.09'
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '''52.03'
$ws.Range('E50').Value = '  -2.78%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '''72.02'
$ws.Range('E51').Value = '  -0.22%  '
